$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for this market/product. It
# belongs chronologically among the existing rows, so insert a fresh row at
# position 179 (this shifts the former rows 179:291 down to 180:292,
# growing the used range from A1:R291 to A1:R292) and fill it with the new
# reading.
$ws.Rows.Item(179).Insert()

$newRow = 179
$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 44981
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112017
$ws.Cells.Item($newRow, 7).Value = "Apio"
$ws.Cells.Item($newRow, 8).Value = "Americana (o)"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 100
$ws.Cells.Item($newRow, 11).Value = 9000
$ws.Cells.Item($newRow, 12).Value = 9500
$ws.Cells.Item($newRow, 13).Value = 9250
$ws.Cells.Item($newRow, 14).Value = "$/docena de matas"
$ws.Cells.Item($newRow, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($newRow, 16).Value = 1542
$ws.Cells.Item($newRow, 17).Value = 6
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
